$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1896103896103896
$ws.Range("C2").Value = 0.548051948051948
$ws.Range("J2").Value = 0.007792207792207792
$ws.Range("P2").Value = 0.1532467532467532
$ws.Range("S2").Value = 0.1012987012987013
$ws.Range("B3").Value = 0.009302325581395349
$ws.Range("C3").Value = 0.02325581395348837
$ws.Range("J3").Value = 0.03720930232558139
$ws.Range("P3").Value = 0.7534883720930232
$ws.Range("S3").Value = 0.1767441860465116
$ws.Range("J4").Value = 0.03508771929824561
$ws.Range("P4").Value = 0.7543859649122807
$ws.Range("S4").Value = 0.2105263157894737
$ws.Range("B6").Value = 0.05
$ws.Range("D6").Value = 0.02727272727272727
$ws.Range("F6").Value = 0.1090909090909091
$ws.Range("J6").Value = 0.2454545454545455
$ws.Range("O6").Value = 0.02727272727272727
$ws.Range("Q6").Value = 0.1181818181818182
$ws.Range("R6").Value = 0.07727272727272727
$ws.Range("S6").Value = 0.3454545454545455
$ws.Range("B7").Value = 0.1395348837209302
$ws.Range("D7").Value = 0.01395348837209302
$ws.Range("F7").Value = 0.03720930232558139
$ws.Range("J7").Value = 0.1255813953488372
$ws.Range("O7").Value = 0.03255813953488372
$ws.Range("Q7").Value = 0.1395348837209302
$ws.Range("R7").Value = 0.07441860465116279
$ws.Range("S7").Value = 0.4372093023255814
$ws.Range("B8").Value = 0.1274725274725275
$ws.Range("D8").Value = 0.02197802197802198
$ws.Range("E8").Value = 0.002197802197802198
$ws.Range("F8").Value = 0.06153846153846154
$ws.Range("J8").Value = 0.1098901098901099
$ws.Range("O8").Value = 0.03736263736263736
$ws.Range("Q8").Value = 0.1692307692307692
$ws.Range("R8").Value = 0.07912087912087912
$ws.Range("S8").Value = 0.3912087912087912
$ws.Range("B9").Value = 0.1589403973509934
$ws.Range("D9").Value = 0.02649006622516556
$ws.Range("F9").Value = 0.04635761589403974
$ws.Range("J9").Value = 0.05960264900662252
$ws.Range("O9").Value = 0.01986754966887417
$ws.Range("Q9").Value = 0.1258278145695364
$ws.Range("R9").Value = 0.0728476821192053
$ws.Range("S9").Value = 0.4900662251655629
$ws.Range("B10").Value = 0.1416275430359937
$ws.Range("D10").Value = 0.02895148669796557
$ws.Range("E10").Value = 0.000782472613458529
$ws.Range("F10").Value = 0.0594679186228482
$ws.Range("J10").Value = 0.1322378716744914
$ws.Range("O10").Value = 0.0242566510172144
$ws.Range("Q10").Value = 0.1705790297339593
$ws.Range("R10").Value = 0.07746478873239436
$ws.Range("S10").Value = 0.3646322378716745
$ws.Range("G11").Value = 0.1358024691358025
$ws.Range("J11").Value = 0.07716049382716049
$ws.Range("K11").Value = 0.1944444444444444
$ws.Range("L11").Value = 0.5709876543209876
$ws.Range("S11").Value = 0.02160493827160494
$ws.Range("G12").Value = 0.7551020408163265
$ws.Range("J12").Value = 0.1836734693877551
$ws.Range("K12").Value = 0.01020408163265306
$ws.Range("L12").Value = 0.04081632653061224
$ws.Range("S12").Value = 0.01020408163265306
$ws.Range("G13").Value = 0.6818181818181818
$ws.Range("J13").Value = 0.2272727272727273
$ws.Range("S13").Value = 0.09090909090909091
$ws.Range("F15").Value = 0.02439024390243903
$ws.Range("H15").Value = 0.1902439024390244
$ws.Range("I15").Value = 0.04390243902439024
$ws.Range("J15").Value = 0.3219512195121951
$ws.Range("K15").Value = 0.05853658536585366
$ws.Range("M15").Value = 0.02439024390243903
$ws.Range("O15").Value = 0.06341463414634146
$ws.Range("S15").Value = 0.2731707317073171
$ws.Range("F16").Value = 0.01574803149606299
$ws.Range("H16").Value = 0.1771653543307087
$ws.Range("I16").Value = 0.06299212598425197
$ws.Range("J16").Value = 0.4448818897637795
$ws.Range("K16").Value = 0.1102362204724409
$ws.Range("M16").Value = 0.007874015748031496
$ws.Range("O16").Value = 0.04724409448818898
$ws.Range("S16").Value = 0.1338582677165354
$ws.Range("F17").Value = 0.03260869565217391
$ws.Range("H17").Value = 0.1684782608695652
$ws.Range("I17").Value = 0.08423913043478261
$ws.Range("J17").Value = 0.4076086956521739
$ws.Range("K17").Value = 0.08695652173913043
$ws.Range("M17").Value = 0.02717391304347826
$ws.Range("N17").Value = 0.002717391304347826
$ws.Range("O17").Value = 0.06793478260869565
$ws.Range("S17").Value = 0.1222826086956522
$ws.Range("F18").Value = 0.005617977528089887
$ws.Range("H18").Value = 0.1853932584269663
$ws.Range("I18").Value = 0.07303370786516854
$ws.Range("J18").Value = 0.4382022471910113
$ws.Range("K18").Value = 0.1292134831460674
$ws.Range("M18").Value = 0.005617977528089887
$ws.Range("O18").Value = 0.05056179775280899
$ws.Range("S18").Value = 0.1123595505617977
$ws.Range("F19").Value = 0.01957713390759593
$ws.Range("H19").Value = 0.2223962411902898
$ws.Range("I19").Value = 0.06499608457321848
$ws.Range("J19").Value = 0.3852779953014879
$ws.Range("K19").Value = 0.1205951448707909
$ws.Range("M19").Value = 0.02192638997650744
$ws.Range("O19").Value = 0.05011746280344558
$ws.Range("S19").Value = 0.1151135473766641
